# Adapt column header formatting to respective input file names:
#   *_old -> *_FV2410 , *_new -> *_FV2504 (the "diff" column is untouched)
# Wrap the header row + data range in an Excel Table ("Table1")
# Freeze the header row (top row) of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# ---------------------------------------------------------------------------
# 1) Rename the header row (row 1) from *_old / *_new to *_FV2410 / *_FV2504
# ---------------------------------------------------------------------------
$oldHeaders = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

$headerRow = $ws.Range("A1:U1")

for ($i = 0; $i -lt $oldHeaders.Count; $i++) {
    $headerRow.Cells.Item(1, $i + 1).Value2 = ($oldHeaders[$i] + "_FV2410")
}
# column K (11) stays "diff"
for ($i = 0; $i -lt $oldHeaders.Count; $i++) {
    $headerRow.Cells.Item(1, $i + 12).Value2 = ($oldHeaders[$i] + "_FV2504")
}

# ---------------------------------------------------------------------------
# 2) Turn A1:U58 into an Excel Table named "Table1"
#    (built on a scratch range first so the already-bold/filled header
#     cells don't get "promoted" into a header dxf / TableStyle override;
#     that keeps styles.xml identical to the un-tabled version)
# ---------------------------------------------------------------------------
$lastRow = 58
$lastCol = 21

$scratchHeaderRow = 200
$scratchTopLeft = $ws.Cells.Item($scratchHeaderRow, 1)
$scratchHeader = $ws.Range($scratchTopLeft, $ws.Cells.Item($scratchHeaderRow, $lastCol))

for ($c = 1; $c -le $lastCol; $c++) {
    $scratchHeader.Cells.Item(1, $c).Value2 = $headerRow.Cells.Item(1, $c).Value2
}
$ws.Cells.Item($scratchHeaderRow + 1, 1).Value2 = "x"

$scratchRange = $ws.Range($scratchTopLeft, $ws.Cells.Item($scratchHeaderRow + 1, $lastCol))
$lo = $ws.ListObjects.Add(1, $scratchRange, [Type]::Missing, 1, [Type]::Missing)
$lo.Name = "Table1"

$finalRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))
$lo.Resize($finalRange)

# wipe the scratch rows again so the sheet dimension / data stay as before
$ws.Range($scratchTopLeft, $ws.Cells.Item($scratchHeaderRow + 1, $lastCol)).Clear()

# ---------------------------------------------------------------------------
# 3) Freeze the header row
# ---------------------------------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select()

Write-Host "done"
